$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 6363.4375
$ws.Range("I11").Value = 6363.4375
$ws.Range("K11").Value = 6363.4375
$ws.Range("M11").Value = -6223.4375

$ws.Range("H32").Value = 1054.0834
$ws.Range("I32").Value = 107.166664
$ws.Range("K32").Value = 107.166664
$ws.Range("M32").Value = 218.833336

$ws.Range("H33").Value = 305.6
$ws.Range("I33").Value = 276.13333
$ws.Range("J33").Value = 394
$ws.Range("K33").Value = 276.13333
$ws.Range("L33").Value = 394
$ws.Range("M33").Value = -47.13333
$ws.Range("N33").Value = -852

$ws.Range("H62").Value = 2489.4
$ws.Range("I62").Value = 2489.4
$ws.Range("K62").Value = 2489.4
$ws.Range("M62").Value = -1865.4

$ws.Range("H65").Value = 2489.4
$ws.Range("I65").Value = 2489.4
$ws.Range("K65").Value = 12447
$ws.Range("M65").Value = -9327

$ws.Range("H74").Value = 6300.615
$ws.Range("I74").Value = 14200
$ws.Range("K74").Value = 14200
$ws.Range("M74").Value = -13264

$ws.Range("H77").Value = 6300.615
$ws.Range("I77").Value = 14200
$ws.Range("K77").Value = 71000
$ws.Range("M77").Value = -66320

$ws.Range("H100").Value = 2939.5833
$ws.Range("I100").Value = 2796.4285
$ws.Range("K100").Value = 2796.4285
$ws.Range("M100").Value = -2255.4285

$ws.Range("H135").Value = 149056.86
$ws.Range("I135").Value = 206465.2
$ws.Range("J135").Value = 5536
$ws.Range("K135").Value = 1858186.8
$ws.Range("L135").Value = 49824
$ws.Range("M135").Value = -1855651.8
$ws.Range("N135").Value = -54894

$ws.Range("H137").Value = 2331009
$ws.Range("I137").Value = 4168732.8
$ws.Range("J137").Value = 9673.842000000001
$ws.Range("K137").Value = 12506198.4
$ws.Range("L137").Value = 29021.526
$ws.Range("M137").Value = -12503648.4
$ws.Range("N137").Value = -34121.526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3270.98
$ws.Range("I32").Value = 3253.5151
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 3253.5151
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -2966.5151
$ws.Range("N32").Value = -5574

$ws.Range("H61").Value = 37112800
$ws.Range("I61").Value = 45501164
$ws.Range("J61").Value = 203989.8
$ws.Range("K61").Value = 45501164
$ws.Range("L61").Value = 203989.8
$ws.Range("M61").Value = -45500952
$ws.Range("N61").Value = -204413.8

$ws.Range("H74").Value = 4595320
$ws.Range("I74").Value = 5225138.5
$ws.Range("K74").Value = 5225138.5
$ws.Range("M74").Value = -5224264.5

$ws.Range("H77").Value = 4595320
$ws.Range("I77").Value = 5225138.5
$ws.Range("K77").Value = 26125692.5
$ws.Range("M77").Value = -26121324.5

$ws.Range("H97").Value = 984.04
$ws.Range("I97").Value = 999.5833
$ws.Range("J97").Value = 611
$ws.Range("K97").Value = 999.5833
$ws.Range("L97").Value = 611
$ws.Range("M97").Value = -503.5833
$ws.Range("N97").Value = -1603

$ws.Range("H136").Value = 37112800
$ws.Range("I136").Value = 45501164
$ws.Range("J136").Value = 203989.8
$ws.Range("K136").Value = 136503492
$ws.Range("L136").Value = 611969.3999999999
$ws.Range("M136").Value = -136500942
$ws.Range("N136").Value = -617069.3999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null

$ws.Range("H94").Value = 927.5
$ws.Range("I94").Value = 902.7273
$ws.Range("K94").Value = 902.7273
$ws.Range("M94").Value = -451.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 42526.72
$ws.Range("I132").Value = 2479.7334
$ws.Range("J132").Value = 102597.2
$ws.Range("K132").Value = 7439.2002
$ws.Range("L132").Value = 307791.6
$ws.Range("M132").Value = -4909.2002
$ws.Range("N132").Value = -312851.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5002158
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 5265429
$ws.Range("K4").Value = 30
$ws.Range("L4").Value = 15796287
$ws.Range("M4").Value = 82
$ws.Range("N4").Value = -15796511

$ws.Range("H107").Value = 745.38635
$ws.Range("I107").Value = 716.12
$ws.Range("J107").Value = 783.8946999999999
$ws.Range("K107").Value = 2148.36
$ws.Range("L107").Value = 2351.6841
$ws.Range("M107").Value = -228.3600000000001
$ws.Range("N107").Value = -6191.6841

$ws.Range("H113").Value = 558.39624
$ws.Range("I113").Value = 507.4091
$ws.Range("J113").Value = 594.5806
$ws.Range("K113").Value = 1522.2273
$ws.Range("L113").Value = 1783.7418
$ws.Range("M113").Value = 647.7727
$ws.Range("N113").Value = -6123.7418

$ws.Range("H120").Value = 4398
$ws.Range("I120").Value = 4398
$ws.Range("K120").Value = 13194
$ws.Range("M120").Value = -8356

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4412.077
$ws.Range("I80").Value = 3879.2
$ws.Range("J80").Value = 4745.125
$ws.Range("K80").Value = 3879.2
$ws.Range("L80").Value = 4745.125
$ws.Range("M80").Value = -2881.2
$ws.Range("N80").Value = -6741.125

$ws.Range("H83").Value = 4412.077
$ws.Range("I83").Value = 3879.2
$ws.Range("J83").Value = 4745.125
$ws.Range("K83").Value = 19396
$ws.Range("L83").Value = 23725.625
$ws.Range("M83").Value = -14404
$ws.Range("N83").Value = -33709.625

$ws.Range("H122").Value = 2166.5
$ws.Range("I122").Value = 2199.0908
$ws.Range("J122").Value = 1808
$ws.Range("K122").Value = 6597.2724
$ws.Range("L122").Value = 5424
$ws.Range("M122").Value = -4147.2724
$ws.Range("N122").Value = -10324

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 800.13336
$ws.Range("I22").Value = 633.3333
$ws.Range("J22").Value = 1050.3334
$ws.Range("K22").Value = 633.3333
$ws.Range("L22").Value = 1050.3334
$ws.Range("M22").Value = -338.3333
$ws.Range("N22").Value = -1640.3334

$ws.Range("H27").Value = 800.13336
$ws.Range("I27").Value = 633.3333
$ws.Range("J27").Value = 1050.3334
$ws.Range("K27").Value = 633.3333
$ws.Range("L27").Value = 1050.3334
$ws.Range("M27").Value = -526.3333
$ws.Range("N27").Value = -1264.3334

$ws.Range("H82").Value = 1047.2222
$ws.Range("I82").Value = 890.625
$ws.Range("K82").Value = 890.625
$ws.Range("M82").Value = -529.625

$ws.Range("H85").Value = 1047.2222
$ws.Range("I85").Value = 890.625
$ws.Range("K85").Value = 890.625
$ws.Range("M85").Value = 357.375

$ws.Range("H93").Value = 1143
$ws.Range("I93").Value = 1097.1904
$ws.Range("K93").Value = 1097.1904
$ws.Range("M93").Value = 150.8096

$ws.Range("H100").Value = 1923.75
$ws.Range("I100").Value = 1821
$ws.Range("K100").Value = 1821
$ws.Range("M100").Value = -1280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 631.3143
$ws.Range("I113").Value = 767
$ws.Range("J113").Value = 427.7857
$ws.Range("K113").Value = 2301
$ws.Range("L113").Value = 1283.3571
$ws.Range("M113").Value = -131
$ws.Range("N113").Value = -5623.3571

$ws.Range("H124").Value = 50233.332
$ws.Range("J124").Value = 50233.332
$ws.Range("L124").Value = 50233.332
$ws.Range("N124").Value = -60053.332

$ws.Range("H135").Value = 54895
$ws.Range("J135").Value = 54895
$ws.Range("L135").Value = 54895
$ws.Range("N135").Value = -65035

$ws.Range("H136").Value = 55254.605
$ws.Range("I136").Value = 39646.54
$ws.Range("K136").Value = 118939.62
$ws.Range("M136").Value = -116389.62
